# Actualización automática 2025-11-28 15:30:09
# Applies November sales updates for asesor HIDALGO HIDALGO PEDRO GUSTAVO
# across the three report sheets (VENTAS POR GRUPO, VENTA MENSUAL,
# CUMPLIMIENTO MENSUAL).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "VENTAS POR GRUPO": per-client sales broken out by product group
# ---------------------------------------------------------------------
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")

# CHASIQUIZA CAMPAÑA JOSE LUIS (row 6): new sales recorded
$wsGrupo.Range("E6").Value = 537.5      # FREGADEROS DE COCINA
$wsGrupo.Range("L6").Value = 1267.08    # PIEDRA SINTERIZADA

# MEGAMAFERS S.A. (row 13): additional PORCELANATO sales
$wsGrupo.Range("M13").Value = 12825.62

# PADILLA MIER BERTHA MARIETA (row 17): additional PORCELANATO sales
$wsGrupo.Range("M17").Value = 4419.18

# Row 23: "N de 21" advisor-reached counters per product group
$wsGrupo.Range("E23").Value = "1 de 21"
$wsGrupo.Range("L23").Value = "3 de 21"

# ---------------------------------------------------------------------
# Sheet "VENTA MENSUAL": per-client sales broken out by month
# ---------------------------------------------------------------------
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")

$wsMensual.Range("F6").Value = 10679.6   # CHASIQUIZA CAMPAÑA JOSE LUIS - noviembre
$wsMensual.Range("F13").Value = 12825.62 # MEGAMAFERS S.A. - noviembre
$wsMensual.Range("F17").Value = 4419.18  # PADILLA MIER BERTHA MARIETA - noviembre
$wsMensual.Range("F23").Value = 63946.43 # TOTAL - noviembre

# ---------------------------------------------------------------------
# Sheet "CUMPLIMIENTO MENSUAL": budget-vs-sales compliance by product group
# ---------------------------------------------------------------------
$wsCumpl = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# FREGADEROS DE COCINA (row 4)
$wsCumpl.Range("D4").Value = 537.5
$wsCumpl.Range("E4").Value = -330.11
$wsCumpl.Range("F4").Value = 2.591735377790636

# PIEDRA SINTERIZADA (row 11)
$wsCumpl.Range("D11").Value = 3604.64
$wsCumpl.Range("E11").Value = -926.6399999999999
$wsCumpl.Range("F11").Value = 1.346019417475728

# PORCELANATO (row 12)
$wsCumpl.Range("D12").Value = 44537.04
$wsCumpl.Range("E12").Value = -119.0400000000009
$wsCumpl.Range("F12").Value = 1.002679994596785

# TOTAL (row 14)
$wsCumpl.Range("D14").Value = 63946.43
$wsCumpl.Range("E14").Value = -8546.958988299059
$wsCumpl.Range("F14").Value = 1.15427871119011
